$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns remain text so numeric-looking strings are not
# auto-converted to numbers by Excel (matches the original inlineStr cell type).
$ws.Columns("D:E").NumberFormat = "@"

$ws.Range('D2').Value = '68.322.59'
$ws.Range('E2').Value = '  -1.13%  '
$ws.Range('D3').Value = '3.900.51'
$ws.Range('E3').Value = '  -1.08%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '486.60'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').Value = '145.88'
$ws.Range('E6').Value = '  -0.41%  '
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '0.741'
$ws.Range('E9').Value = '  +2.36%  '
$ws.Range('D10').Value = '0.182'
$ws.Range('E10').Value = '  +6.90%  '
$ws.Range('E11').Value = '  -0.52%  '
$ws.Range('D12').Value = '42.93'
$ws.Range('E12').Value = '  +0.55%  '
$ws.Range('D13').Value = '10.47'
$ws.Range('E13').Value = '  +0.00%  '
$ws.Range('D14').Value = '4.521.91'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '3.935.38'
$ws.Range('E15').Value = '  -1.00%  '
$ws.Range('D16').Value = '14.19'
$ws.Range('E16').Value = '  -2.85%  '
$ws.Range('E17').Value = '  -0.67%  '
$ws.Range('D18').Value = '20.00'
$ws.Range('E18').Value = '  +1.42%  '
$ws.Range('E19').Value = '  +0.41%  '
$ws.Range('D20').Value = '68.339.97'
$ws.Range('D21').Value = '431.14'
$ws.Range('E21').Value = '  -1.00%  '
$ws.Range('E22').Value = '  +6.36%  '
$ws.Range('D23').Value = '14.76'
$ws.Range('E23').Value = '  +1.07%  '
$ws.Range('D24').Value = '12.48'
$ws.Range('E24').Value = '  +21.70%  '
$ws.Range('E25').Value = '  +1.39%  '
$ws.Range('D26').Value = '3.71'
$ws.Range('E26').Value = '  +3.84%  '
$ws.Range('D27').Value = '10.96'
$ws.Range('E27').Value = '  -6.27%  '
$ws.Range('D28').Value = '37.30'
$ws.Range('E28').Value = '  -2.84%  '
$ws.Range('E29').Value = '  -3.57%  '
$ws.Range('D30').Value = '721.08'
$ws.Range('E30').Value = '  +1.02%  '
$ws.Range('D31').Value = '13.42'
$ws.Range('E31').Value = '  +1.03%  '
$ws.Range('E32').Value = '  +0.29%  '
$ws.Range('E33').Value = '  +2.45%  '
$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').Value = '0.0₃0889'
$ws.Range('E34').Value = '  -4.69%  '
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '61.64'
$ws.Range('E35').Value = '  +4.87%  '
$ws.Range('E36').Value = '  +7.31%  '
$ws.Range('D37').Value = '40.80'
$ws.Range('E37').Value = '  -1.28%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').Value = '0.147'
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('B39').Value = 'TheGraph'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D39').Value = '0.400'
$ws.Range('E39').Value = '  +17.00%  '
$ws.Range('D40').Value = '0.997'
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '3.00'
$ws.Range('E41').Value = '  +8.75%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = '0.0495'
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('D43').Value = '3.07'
$ws.Range('E43').Value = '  +3.00%  '
$ws.Range('D44').Value = '3.01'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').Value = '0.0₆0373'
$ws.Range('E45').Value = '  +31.92%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '0.142'
$ws.Range('E46').Value = '  +0.63%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D47').Value = '3.34'
$ws.Range('E47').Value = '  +7.18%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.02%  '
$ws.Range('E49').Value = '  -1.08%  '
$ws.Range('E50').Value = '  -2.68%  '
$ws.Range('D51').Value = '144.20'
$ws.Range('E51').Value = '  -2.57%  '
